# Add a new "Save" column (H) to the s_vals sheet, mirroring the header
# style used by the other stat columns (B:G) and filling the two data
# rows with 0, matching the committed diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) from the neighboring "sum" header (G1) onto
# the new H1 header cell so it reuses the same bold/border/centered xf
# instead of minting a brand-new style.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Header text for the new column.
$ws.Range("H1").Value = "Save"

# Data rows: both existing rows get a 0 in the new column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
